$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old (E5, F5) values that are being replaced
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()

# Write the new data block C5:D10
$data = @(
    @(0.0199743, 0.0026535),
    @(0.0177186, 0.0048243),
    @(0.0057744, 0.0062909),
    @(0.0188013, 0.0033218),
    @(0.019903,  0.0026143),
    @(0.0199379, 0.0042961)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}
